$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the interactionPath column (E) to use a generic templated value
# for every data row, since the actual path is now auto-generated from
# dayID and interactionID rather than being entered manually per row.
$ws.Range("E2:E5").Value = "Cashier_Interaction_Day##_**"

# Add explanatory notes in column I describing the template placeholders.
$ws.Range("I4").Value = "** is interactionID"
$ws.Range("I3").Value = "## is dayID"
$ws.Range("I2").Value = "AUTO"

# Set the new column width for column I to match the authored width
# (~20.33 characters wide, i.e. "20.33203125" in the saved OOXML).
$ws.Columns.Item(9).ColumnWidth = 19.5

# Update the selection to match what was recorded when the file was saved.
$ws.Range("I3").Select()

$wb.Save()
